# Add columns I (I0) and J (IF) to Sheet1, mirroring headers/styles of
# existing column H, and fill in the per-row values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header style used by the existing header row (e.g. H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for rows 2-32 (columns I and J)
$iValues = @(9,5,8,3,7,8,7,9,7,6,7,8,7,8,9,9,8,8,9,8,8,6,8,7,9,8,7,6,8,7,6)
$jValues = @(9,6,8,4,7,8,7,9,7,6,7,8,7,8,9,9,8,8,9,8,8,6,8,8,9,8,8,6,8,7,6)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
